# Financials update: insert a new "most recent period" column before column D,
# shifting the existing D:K data right to E:L, then populate the new column D
# with the latest period's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new blank column at D; everything in D:K shifts to E:L.
$ws.Columns("D:D").Insert()

# 2) Copy the number formats from the (now shifted) old column, column E,
#    into the freshly inserted column D so the new cells carry the same
#    date / number styling as the rest of each row.
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)

# 3) Populate the new column D with the latest period's values.
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 279332000
$ws.Range("D9").Value = 192854000
$ws.Range("D10").Value = 86478000
$ws.Range("D12").Value = 1466000
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = "NA"
$ws.Range("D15").Value = 18745000
$ws.Range("D17").Value = 258493000
$ws.Range("D18").Value = 20839000
$ws.Range("D20").Value = 10880000
$ws.Range("D21").Value = 50464000
$ws.Range("D22").Value = 766000
$ws.Range("D23").Value = 30953000
$ws.Range("D24").Value = 9823000
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 21130000
$ws.Range("D27").Value = 20549000
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 291000
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -10880000
$ws.Range("D33").Value = 20840000
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 20840000
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 3042000
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 24701000
$ws.Range("D44").Value = 18958000
$ws.Range("D45").Value = 1272000
$ws.Range("D46").Value = 47973000
$ws.Range("D47").Value = 32182000
$ws.Range("D48").Value = 247101000
$ws.Range("D49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 18940000
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 346196000
$ws.Range("D57").Value = 21063000
$ws.Range("D58").Value = 17258000
$ws.Range("D59").Value = 18817000
$ws.Range("D60").Value = 57138000
$ws.Range("D61").Value = 20538000
$ws.Range("D62").Value = 69992000
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 154402000
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 421653000
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 191794000
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 20840000
$ws.Range("D83").Value = 18745000
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 36014000
$ws.Range("D91").Value = -19574000
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -16446000
$ws.Range("D96").Value = -13798000
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -19446000
$ws.Range("D101").Value = -257000
$ws.Range("D102").Value = -135000
